# Update DM integration fixture hierarchies
# Updates the ID (UUID) columns on several sheets with newly generated
# identifiers, and widens column A on the CodeSchemes, Extensions and
# Members_dpmTypedDomain sheets to fit the new values.

$wb = $excel.ActiveWorkbook

# --- CodeSchemes sheet -------------------------------------------------
$ws = $wb.Worksheets.Item("CodeSchemes")
$ws.Columns.Item(1).ColumnWidth = 35.57142857142857
$ws.Range("A2").Value = "41428e0e-16e0-4761-bf43-3a1221e125d4"

# --- Codes sheet ---------------------------------------------------------
$ws = $wb.Worksheets.Item("Codes")
$ws.Range("A2").Value = "eb2c76c7-3e52-4afa-9ee9-bf06d0c801c4"
$ws.Range("A3").Value = "701c21b8-3b02-4adc-9acb-79f6be5fbe8d"
$ws.Range("A4").Value = "eaefdc8e-45a2-41ec-9cc4-99d28d687ef4"
$ws.Range("A5").Value = "7b9c3147-a9d6-4c96-ac35-f741651e721e"
$ws.Range("A6").Value = "32d88269-a59e-4505-92ac-7fd2a7e053e7"
$ws.Range("A7").Value = "d3a2e1fe-a9ba-49e8-a553-5b41094aa138"
$ws.Range("A8").Value = "fe5f0c87-5b80-4e9c-a0a3-6e34814c0fde"
$ws.Range("A9").Value = "6ab3e74d-285a-4e22-887e-681bc696a1cb"

# --- Extensions sheet ------------------------------------------------------
$ws = $wb.Worksheets.Item("Extensions")
$ws.Columns.Item(1).ColumnWidth = 32.285714285714285
$ws.Range("A2").Value = "8134f246-8059-4435-a4d2-06cc1ce088e1"

# --- Members_dpmTypedDomain sheet ------------------------------------------
$ws = $wb.Worksheets.Item("Members_dpmTypedDomain")
$ws.Columns.Item(1).ColumnWidth = 35.57142857142857
$ws.Range("A2").Value = "5791e8c1-a618-4c45-97f6-f67a46bdbf96"
$ws.Range("A3").Value = "78a55fe2-6b50-4274-a1cb-3d8e2ff7bea9"
$ws.Range("A4").Value = "8dd2f75d-ffb9-4a99-964a-230287adbdf6"
$ws.Range("A5").Value = "9e126408-7159-4871-b5ab-7bb0343eda93"
$ws.Range("A6").Value = "c8f8f88b-6ae1-4efd-bb34-35e33b675e0a"
$ws.Range("A7").Value = "a5392778-b3d5-4bdc-bda7-74b19aa92d72"
$ws.Range("A8").Value = "f416fad4-24cd-41b2-97d9-c076b4e4f2e3"
$ws.Range("A9").Value = "961c73d4-89b6-4513-9eb7-aa7719a236f6"
